# Fit the Weibull distribution
# - Remove the "Industrial chemical" sheet entirely.
# - Update the model-summary coefficients (Value / Std. Error / p) on the
#   remaining chemical-category sheets to the refit (Weibull) numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Drop the "Industrial chemical" worksheet.
# ---------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Industrial chemical")
$wsDelete.Delete()

# ---------------------------------------------------------------------
# 2) Helper: push an array of [row, B, C, D] rows into a sheet's B:D
#    columns (rows 2-24).
# ---------------------------------------------------------------------
function Set-ModelValues($sheetName, $rows) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $ws.Cells.Item($r[0], 2).Value = $r[1]
        $ws.Cells.Item($r[0], 3).Value = $r[2]
        $ws.Cells.Item($r[0], 4).Value = $r[3]
    }
}

# ---------------------------------------------------------------------
# API
# ---------------------------------------------------------------------
$apiRows = @(
    @(2, 1.0719, 6147.9054, 0.9999),
    @(3, -0.3268, 0.2382, 0.17),
    @(4, -0.1697, 0.1871, 0.3644),
    @(5, 1.0187, 6147.9054, 0.9999),
    @(6, 1.1801, 6147.9054, 0.9998),
    @(7, -0.1754, 6160.976, 1),
    @(8, -0.1637, 23587.6925, 1),
    @(9, 0.5556, 6147.9054, 0.9999),
    @(10, -0.22, 6177.3998, 1),
    @(11, -0.0765, 16756.1251, 1),
    @(12, -0.6316, 0.654, 0.3341),
    @(13, -1.1726, 0.8919, 0.1886),
    @(14, -1.3815, 0.9042, 0.1265),
    @(15, -1.3359, 0.8403, 0.1119),
    @(16, -1.1726, 0.7894, 0.1374),
    @(17, -0.91, 0.7453, 0.2221),
    @(18, -0.8464, 0.7287, 0.2454),
    @(19, -1.1697, 0.7287, 0.1085),
    @(20, -1.3015, 0.7867, 0.098),
    @(21, -1.2545, 0.8756, 0.1519),
    @(22, -1.1691, 1.171, 0.3181),
    @(23, -1.0838, 1.7815, 0.543),
    @(24, -0.7614, 0.2437, 0.0018)
)
Set-ModelValues "API" $apiRows

# ---------------------------------------------------------------------
# Anthropogenic pollution
# ---------------------------------------------------------------------
$anthroRows = @(
    @(2, 3.0173, 0.3728, 0),
    @(3, 0.19, 0.1014, 0.061),
    @(4, 0.0258, 0.0949, 0.786),
    @(5, -0.3186, 0.1797, 0.0762),
    @(6, 0.4603, 0.164, 0.005),
    @(7, -0.6218, 0.2105, 0.0031),
    @(8, 0.3614, 0.1707, 0.0343),
    @(9, -0.4868, 0.2259, 0.0312),
    @(10, 0.3827, 0.1594, 0.0164),
    @(11, -0.1948, 0.1808, 0.2812),
    @(12, -0.3746, 0.2636, 0.1553),
    @(13, -0.6891, 0.3711, 0.0633),
    @(14, -0.837, 0.3889, 0.0314),
    @(15, -0.9047, 0.388, 0.0197),
    @(16, -0.8325, 0.3752, 0.0265),
    @(17, -0.7239, 0.3594, 0.044),
    @(18, -0.6939, 0.3492, 0.0469),
    @(19, -0.7323, 0.3497, 0.0363),
    @(20, -0.6578, 0.3608, 0.0682),
    @(21, -0.6246, 0.3637, 0.0859),
    @(22, -0.8915, 0.3769, 0.018),
    @(23, -1.2286, 0.4868, 0.0116),
    @(24, -0.3838, 0.0526, 0)
)
Set-ModelValues "Anthropogenic pollution" $anthroRows

# ---------------------------------------------------------------------
# PAH
# ---------------------------------------------------------------------
$pahRows = @(
    @(2, 3.0154, 0.2362, 0),
    @(3, 0.0007, 0.0633, 0.9909),
    @(4, 0.049, 0.0582, 0.4003),
    @(5, 0.0258, 0.1063, 0.8078),
    @(6, 0.0984, 0.1105, 0.3729),
    @(7, 0.0582, 0.1088, 0.5927),
    @(8, -0.0619, 0.1145, 0.5889),
    @(9, -0.2267, 0.1372, 0.0985),
    @(10, -0.1448, 0.111, 0.1922),
    @(11, -0.1693, 0.119, 0.1547),
    @(12, -0.173, 0.1649, 0.294),
    @(13, -0.3233, 0.2305, 0.1607),
    @(14, -0.3871, 0.238, 0.1039),
    @(15, -0.4319, 0.2329, 0.0637),
    @(16, -0.4815, 0.2254, 0.0327),
    @(17, -0.5076, 0.2194, 0.0207),
    @(18, -0.475, 0.2142, 0.0266),
    @(19, -0.4041, 0.2133, 0.0582),
    @(20, -0.3984, 0.2221, 0.0728),
    @(21, -0.4137, 0.2277, 0.0693),
    @(22, -0.3643, 0.2334, 0.1185),
    @(23, -0.2967, 0.3006, 0.3236),
    @(24, -0.9546, 0.0661, 0)
)
Set-ModelValues "PAH" $pahRows

# ---------------------------------------------------------------------
# PCP
# ---------------------------------------------------------------------
$pcpRows = @(
    @(2, 3.7515, 0.4089, 0),
    @(3, 0.0611, 0.0838, 0.4659),
    @(4, 0.0081, 0.0822, 0.9215),
    @(5, 0.4339, 0.2323, 0.0618),
    @(6, 0.1996, 0.2332, 0.392),
    @(7, 0.2214, 0.2435, 0.3633),
    @(8, -0.0493, 0.3057, 0.8719),
    @(9, 0.0826, 0.296, 0.7803),
    @(10, -0.0255, 0.2845, 0.9286),
    @(11, 0.1265, 0.2355, 0.5912),
    @(12, -0.2965, 0.2866, 0.3009),
    @(13, -0.5605, 0.4005, 0.1617),
    @(14, -0.6597, 0.4219, 0.1179),
    @(15, -0.6071, 0.4086, 0.1373),
    @(16, -0.4944, 0.3736, 0.1857),
    @(17, -0.3852, 0.3505, 0.2718),
    @(18, -0.2984, 0.3429, 0.3842),
    @(19, -0.2582, 0.3434, 0.4521),
    @(20, -0.2897, 0.3607, 0.4218),
    @(21, -0.3663, 0.3838, 0.3399),
    @(22, -0.4244, 0.4507, 0.3464),
    @(23, -0.4789, 0.6402, 0.4544),
    @(24, -1.0572, 0.1036, 0)
)
Set-ModelValues "PCP" $pcpRows

# ---------------------------------------------------------------------
# POP
# ---------------------------------------------------------------------
$popRows = @(
    @(2, 2.2819, 0.3198, 0),
    @(3, -0.3976, 0.0877, 0),
    @(4, -0.4152, 0.091, 0),
    @(5, 0.4098, 0.1705, 0.0162),
    @(6, -0.3434, 0.2282, 0.1324),
    @(7, 0.2732, 0.1781, 0.1251),
    @(8, -0.4675, 0.2682, 0.0813),
    @(9, 0.4277, 0.1958, 0.029),
    @(10, 0.1174, 0.1772, 0.5076),
    @(11, -0.4106, 0.2805, 0.1432),
    @(12, -0.0623, 0.2198, 0.7768),
    @(13, -0.1169, 0.3014, 0.6981),
    @(14, -0.1194, 0.3122, 0.7022),
    @(15, -0.0318, 0.3081, 0.9177),
    @(16, 0.0892, 0.2944, 0.7619),
    @(17, 0.2285, 0.2871, 0.4261),
    @(18, 0.2597, 0.2835, 0.3596),
    @(19, 0.048, 0.2842, 0.8658),
    @(20, -0.0453, 0.2999, 0.88),
    @(21, 0.0083, 0.3088, 0.9786),
    @(22, 0.1138, 0.3171, 0.7196),
    @(23, 0.2256, 0.4173, 0.5887),
    @(24, -0.9581, 0.0612, 0)
)
Set-ModelValues "POP" $popRows

# ---------------------------------------------------------------------
# Pesticide
# ---------------------------------------------------------------------
$pesticideRows = @(
    @(2, 3.4383, 1.18, 0.0036),
    @(3, 0.6208, 0.2413, 0.0101),
    @(4, 0.3603, 0.2329, 0.1218),
    @(5, -0.0803, 0.5219, 0.8778),
    @(6, 0.5088, 0.3949, 0.1975),
    @(7, -0.4649, 0.6191, 0.4527),
    @(8, 0.4744, 0.3789, 0.2106),
    @(9, -0.6808, 0.6625, 0.3041),
    @(10, -1.6174, 942.2577, 0.9986),
    @(11, -1.0266, 1300.57, 0.9994),
    @(12, -1.1952, 0.9579, 0.2121),
    @(13, -2.2438, 1.4255, 0.1155),
    @(14, -2.8393, 1.6632, 0.0878),
    @(15, -3.0455, 1.6886, 0.0713),
    @(16, -2.9566, 1.5161, 0.0512),
    @(17, -2.6644, 1.2651, 0.0352),
    @(18, -2.1991, 1.1244, 0.0505),
    @(19, -2.1341, 1.1177, 0.0562),
    @(20, -2.2101, 1.1493, 0.0545),
    @(21, -2.0527, 1.1667, 0.0785),
    @(22, -2.2439, 1.2008, 0.0617),
    @(23, -2.505, 1.5055, 0.0961),
    @(24, -0.6015, 0.1774, 0.0007)
)
Set-ModelValues "Pesticide" $pesticideRows

# ---------------------------------------------------------------------
# Plasticizer
# ---------------------------------------------------------------------
$plasticizerRows = @(
    @(2, 3.3102, 0.6893, 0),
    @(3, -0.0768, 0.2004, 0.7013),
    @(4, 0.3973, 0.1766, 0.0245),
    @(5, 0.0691, 0.3382, 0.8382),
    @(6, -0.0892, 0.3434, 0.7949),
    @(7, 0.1367, 0.3414, 0.6888),
    @(8, 0.1667, 0.3501, 0.6339),
    @(9, -0.282, 0.4135, 0.4952),
    @(10, -0.8794, 0.3386, 0.0094),
    @(11, 0.432, 0.3456, 0.2113),
    @(12, -0.3069, 0.4892, 0.5305),
    @(13, -0.541, 0.6848, 0.4295),
    @(14, -0.4909, 0.7054, 0.4865),
    @(15, -0.4449, 0.6915, 0.52),
    @(16, -0.5284, 0.6708, 0.4309),
    @(17, -0.6134, 0.6539, 0.3483),
    @(18, -0.6791, 0.6436, 0.2913),
    @(19, -0.7422, 0.6429, 0.2483),
    @(20, -0.8439, 0.6685, 0.2068),
    @(21, -0.7499, 0.6845, 0.2733),
    @(22, -0.5143, 0.7125, 0.4704),
    @(23, -0.2765, 0.9239, 0.7647),
    @(24, 0.2522, 0.0658, 0.0001)
)
Set-ModelValues "Plasticizer" $plasticizerRows
